# Update "想去人数" (want-to-go count) figures in column F across the
# 展览 / 本地生活 / 全部类型 sheets to match the latest scrape output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 323
$ws1.Range("F14").Value = 922
$ws1.Range("F22").Value = 3279
$ws1.Range("F23").Value = 967
$ws1.Range("F25").Value = 2356
$ws1.Range("F28").Value = 3236
$ws1.Range("F32").Value = 1422
$ws1.Range("F34").Value = 761
$ws1.Range("F37").Value = 93
$ws1.Range("F39").Value = 1165
$ws1.Range("F40").Value = 1840
$ws1.Range("F41").Value = 433
$ws1.Range("F47").Value = 61

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 158

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 323
$ws4.Range("F9").Value  = 922
$ws4.Range("F18").Value = 3279
$ws4.Range("F19").Value = 967
$ws4.Range("F22").Value = 2356
$ws4.Range("F24").Value = 3236
$ws4.Range("F34").Value = 761
$ws4.Range("F41").Value = 1165
$ws4.Range("F42").Value = 1840
$ws4.Range("F44").Value = 433
$ws4.Range("F49").Value = 61
